$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (HOUR 1)
$ws.Range("B2").Value = 12500
$ws.Range("C2").Value = 10000
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 22500

# Row 3 (HOUR 2)
$ws.Range("B3").Value = 12500
$ws.Range("C3").Value = 10000
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 22500

# Row 4 (HOUR 3)
$ws.Range("B4").Value = 12500
$ws.Range("C4").Value = 10000
$ws.Range("E4").Value = 22500

# Row 5 (HOUR 4)
$ws.Range("B5").Value = 12500
$ws.Range("C5").Value = 10000
$ws.Range("E5").Value = 22500

# Row 6 (HOUR 5)
$ws.Range("C6").Value = 10000
$ws.Range("E6").Value = 22500

# Row 10 (HOUR 9)
$ws.Range("B10").Value = 25000
$ws.Range("C10").Value = 20000
$ws.Range("D10").Value = 12000
$ws.Range("E10").Value = 57000

# Row 11 (HOUR 10)
$ws.Range("B11").Value = 25000
$ws.Range("D11").Value = 20000
$ws.Range("E11").Value = 65000

# Row 12 (HOUR 11)
$ws.Range("D12").Value = 20000
$ws.Range("E12").Value = 65000
